$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.744.68'
$ws.Range('E2').Value = '  -2.39%  '
$ws.Range('D3').Value = '3.144.08'
$ws.Range('E3').Value = '  -8.05%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '564.66'
$ws.Range('E5').Value = '  -3.32%  '
$ws.Range('D6').Value = '170.69'
$ws.Range('E6').Value = '  -4.53%  '
$ws.Range('D7').Value = '0.617'
$ws.Range('E7').Value = '  -0.79%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '3.142.19'
$ws.Range('E9').Value = '  -8.10%  '
$ws.Range('E10').Value = '  -5.62%  '
$ws.Range('E11').Value = '  -5.78%  '
$ws.Range('D12').Value = '0.394'
$ws.Range('E12').Value = '  -4.73%  '
$ws.Range('D13').Value = '3.688.61'
$ws.Range('E13').Value = '  -8.21%  '
$ws.Range('D14').Value = '0.135'
$ws.Range('E14').Value = '  +0.87%  '
$ws.Range('D15').Value = '27.07'
$ws.Range('E15').Value = '  -8.04%  '
$ws.Range('D16').Value = '64.681.70'
$ws.Range('E16').Value = '  -2.63%  '
$ws.Range('E17').Value = '  -5.81%  '
$ws.Range('D18').Value = '3.145.23'
$ws.Range('E18').Value = '  -8.21%  '
$ws.Range('D19').Value = '5.70'
$ws.Range('E19').Value = '  -3.52%  '
$ws.Range('E20').Value = '  -6.73%  '
$ws.Range('D21').Value = '354.92'
$ws.Range('E21').Value = '  -3.23%  '
$ws.Range('E22').Value = '  -4.85%  '
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.46%  '
$ws.Range('E24').Value = '  -5.97%  '
$ws.Range('E25').Value = '  -6.89%  '
$ws.Range('E26').Value = '  -6.89%  '
$ws.Range('D27').Value = '9.61'
$ws.Range('E27').Value = '  -2.18%  '
$ws.Range('E28').Value = '  -2.36%  '
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.15%  '
$ws.Range('E30').Value = '  -0.15%  '
$ws.Range('E31').Value = '  -4.90%  '
$ws.Range('E32').Value = '  -7.01%  '
$ws.Range('D33').Value = '21.97'
$ws.Range('E33').Value = '  -6.22%  '
$ws.Range('E34').Value = '  -5.37%  '
$ws.Range('E35').Value = '  -5.02%  '
$ws.Range('E36').Value = '  -7.00%  '
$ws.Range('D37').Value = '153.37'
$ws.Range('E37').Value = '  -5.88%  '
$ws.Range('D38').Value = '0.830'
$ws.Range('E38').Value = '  -5.19%  '
$ws.Range('D39').Value = '26.00'
$ws.Range('E39').Value = '  -6.20%  '
$ws.Range('D40').Value = '1.74'
$ws.Range('E40').Value = '  -3.12%  '
$ws.Range('D41').Value = '2.54'
$ws.Range('E41').Value = '  -1.70%  '
$ws.Range('D42').Value = '2.653.62'
$ws.Range('E42').Value = '  -1.90%  '
$ws.Range('D43').Value = '4.18'
$ws.Range('E43').Value = '  -6.25%  '
$ws.Range('E44').Value = '  -5.12%  '
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D45').Value = '0.0654'
$ws.Range('E45').Value = '  -5.19%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').Value = '24.13'
$ws.Range('E46').Value = '  -4.58%  '
$ws.Range('E47').Value = '  -2.38%  '
$ws.Range('D48').Value = '319.97'
$ws.Range('E48').Value = '  -4.08%  '
$ws.Range('D49').Value = '0.0273'
$ws.Range('E49').Value = '  -4.38%  '
$ws.Range('E50').Value = '  -2.35%  '
$ws.Range('E51').Value = '  -0.09%  '
